# normalisasi.xlsx update:
#  - fix " rs " / "korsel" normalisation entries
#  - drop the " vs " / "versus" pair
#  - add hezbollah/hizbullah, -as /-ri , palestine/palestina, kurbo/kubro,
#    pengungsian/pengungsi normalisation rows
#  - add "stemming" (asia/asian) and "NER" (israel/palestina/gaza) sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # normalisasi
$ws2 = $wb.Worksheets.Item(2)   # remove

# ---------------------------------------------------------------------------
# 1. normalisasi: drop the " vs " / "versus" row (old row 25)
# ---------------------------------------------------------------------------
$ws1.Rows.Item(25).Delete()

# ---------------------------------------------------------------------------
# 2. normalisasi: append the new rows (order chosen to control shared-string
#    table layout, matching how they were authored)
# ---------------------------------------------------------------------------
$ws1.Range("A33").Value = "hezbollah"
$ws1.Range("B33").Value = "hizbullah"

$ws1.Range("A34").Value = "-as "
$ws1.Range("B34").Value = "amerika serikat"

$ws1.Range("A35").Value = "-ri "
$ws1.Range("B35").Value = "indonesia"

# fix "rs" -> " rs "
$ws1.Range("A11").Value = " rs "

$ws1.Range("A36").Value = "palestine"
$ws1.Range("B36").Value = "palestina"

# fix "koresel" -> "korsel"
$ws1.Range("A20").Value = "korsel"

# fix stray capitalised "Indonesia" -> "indonesia" (row shifted up after the
# " vs " row was removed)
$ws1.Range("B28").Value = "indonesia"

# ---------------------------------------------------------------------------
# 3. add the "stemming" sheet (asia/asian), placed right after "remove"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "stemming"
$ws3.Range("A1").Value = "asia"
$ws3.Range("B1").Value = "asia"
$ws3.Range("A2").Value = "asian"
$ws3.Range("B2").Value = "asia"

# ---------------------------------------------------------------------------
# 4. normalisasi: remaining new rows (kurbo/kubro, pengungsian/pengungsi)
# ---------------------------------------------------------------------------
$ws1.Range("A37").Value = "kurbo"
$ws1.Range("B37").Value = "kubro"

$ws1.Range("A38").Value = "pengungsian"
$ws1.Range("B38").Value = "pengungsi"

# ---------------------------------------------------------------------------
# 5. add the "NER" sheet (israel/palestina/gaza), placed after "stemming"
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet2)
$ws4.Name = "NER"
$ws4.Range("A1").Value = "israel"
$ws4.Range("A2").Value = "palestina"
$ws4.Range("A3").Value = "gaza"

# ---------------------------------------------------------------------------
# 6. view/selection tweaks
# ---------------------------------------------------------------------------
$ws1.Range("A17").Select()
$ws1.Range("C39").Select()

$ws2.Range("N11").Select()

$ws3.Range("B3").Select()

$ws4.Activate()
$ws4.Range("F25").Select()
